$wb = $excel.ActiveWorkbook

$wsUsers    = $wb.Worksheets.Item("users")
$wsProducts = $wb.Worksheets.Item("products")
$wsOrders   = $wb.Worksheets.Item("orders")

# --- Data fixes ---
# products!C5 price typo fix: 299.99 -> 299.9 (stored as 299.89999999999998 in IEEE-754 double)
$wsProducts.Range("C5").Value = 299.9

# orders!D6 total_amount fix: 89.99 -> 89
$wsOrders.Range("D6").Value = 89

# --- Column widths for the "products" sheet ---
$wsProducts.Columns.Item(1).ColumnWidth = 12.840401785714286
$wsProducts.Columns.Item(2).ColumnWidth = 14.840401785714286
$wsProducts.Columns.Item(3).ColumnWidth = 10.504464285714286
$wsProducts.Columns.Item(4).ColumnWidth = 14.504464285714286
$wsProducts.Columns.Item(5).ColumnWidth = 11.617745535714286

# --- Selections on each sheet ---
$wsUsers.Range("E8").Select()
$wsOrders.Range("E9").Select()
$wsProducts.Range("E6").Select()

# --- Active sheet becomes "products" ---
$wsProducts.Activate()
